$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need to be forced to Text
# so Excel stores them as strings (matching the source inlineStr cells) rather
# than coercing them into numeric cells.
$textCells = @('D6', 'D8', 'D9', 'D10', 'D14', 'D15', 'D16', 'D18', 'D20', 'D21', 'D23', 'D24', 'D25', 'D26', 'D27', 'D29', 'D30', 'D31', 'D32', 'D33', 'D36', 'D37', 'D39', 'D41', 'D42', 'D44', 'D45', 'D48', 'D49', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '29.389.27'
$ws.Range('E2').Value = '  -0.02%  '
$ws.Range('D3').Value = '1.846.34'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('D6').Value = '0.6296'
$ws.Range('E6').Value = '  -1.01%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '0.07591'
$ws.Range('E8').Value = '  +0.46%  '
$ws.Range('D9').Value = '0.2929'
$ws.Range('E9').Value = '  -1.36%  '
$ws.Range('D10').Value = '24.47'
$ws.Range('E11').Value = '  +0.00%  '
$ws.Range('D12').Value = '1.841.80'
$ws.Range('E12').Value = '  +0.09%  '
$ws.Range('E13').Value = '  +9.52%  '
$ws.Range('D14').Value = '5.001'
$ws.Range('D15').Value = '0.6777'
$ws.Range('E15').Value = '  -0.90%  '
$ws.Range('D16').Value = '83.63'
$ws.Range('E16').Value = '  +0.71%  '
$ws.Range('D17').Value = '2.093.52'
$ws.Range('E17').Value = '  -7.54%  '
$ws.Range('D18').Value = '6.143'
$ws.Range('E18').Value = '  -0.55%  '
$ws.Range('D19').Value = '29.408.66'
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('D20').Value = '228.46'
$ws.Range('E20').Value = '  -0.46%  '
$ws.Range('D21').Value = '12.43'
$ws.Range('E21').Value = '  -0.29%  '
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').Value = '7.417'
$ws.Range('E23').Value = '  -1.97%  '
$ws.Range('D24').Value = '1.001'
$ws.Range('E24').Value = '  +0.05%  '
$ws.Range('D25').Value = '157.43'
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('D26').Value = '0.1394'
$ws.Range('E26').Value = '  -0.65%  '
$ws.Range('D27').Value = '8.375'
$ws.Range('E27').Value = '  -0.22%  '
$ws.Range('E28').Value = '  -0.35%  '
$ws.Range('D29').Value = '1.461'
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('D30').Value = '1.301'
$ws.Range('E30').Value = '  +4.10%  '
$ws.Range('D31').Value = '0.05604'
$ws.Range('E31').Value = '  -1.75%  '
$ws.Range('D32').Value = '4.101'
$ws.Range('E32').Value = '  -0.62%  '
$ws.Range('D33').Value = '4.029'
$ws.Range('E33').Value = '  -0.10%  '
$ws.Range('E34').Value = '  -0.40%  '
$ws.Range('E35').Value = '  -0.12%  '
$ws.Range('D36').Value = '0.7101'
$ws.Range('E36').Value = '  -0.82%  '
$ws.Range('D37').Value = '2.585'
$ws.Range('E37').Value = '  -0.41%  '
$ws.Range('D38').Value = '1.229.58'
$ws.Range('E38').Value = '  -1.82%  '
$ws.Range('D39').Value = '2.769'
$ws.Range('E39').Value = '  -0.58%  '
$ws.Range('E40').Value = '  -0.52%  '
$ws.Range('D41').Value = '6.445'
$ws.Range('E41').Value = '  +3.96%  '
$ws.Range('D42').Value = '0.9068'
$ws.Range('E42').Value = '  -0.18%  '
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('D44').Value = '101.65'
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').Value = '65.96'
$ws.Range('E45').Value = '  -0.87%  '
$ws.Range('E46').Value = '  +3.14%  '
$ws.Range('E47').Value = '  +1.41%  '
$ws.Range('D48').Value = '0.4016'
$ws.Range('E48').Value = '  -0.22%  '
$ws.Range('D49').Value = '8.963'
$ws.Range('E49').Value = '  -2.07%  '
$ws.Range('E50').Value = '  -1.41%  '
$ws.Range('D51').Value = '0.1121'
$ws.Range('E51').Value = '  -0.74%  '

# Strip the temporary Text number-format again so the cells end up with no
# explicit style, matching the original workbook formatting.
foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}
